# Applies:
#  1) slide 5's table style id change
#     {D9AA77B7-28C2-44C8-9CAE-CD422E885724} -> {76FC8E9B-0E7A-46F4-B581-E791D2C43184}
#  2) the deck's active theme (Integral / "Red Violet" scheme) is swapped back to the
#     stock "Office Theme" / "Office" colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $sh = $s5.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{76FC8E9B-0E7A-46F4-B581-E791D2C43184}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
# Office Theme ("Office" colour scheme) RGB values, in MsoThemeColorSchemeIndex
# order (Dark1, Light1, Dark2, Light2, Accent1..6, Hyperlink, FollowedHyperlink).
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    $tcs.Colors($i).RGB = $r -bor ($g -shl 8) -bor ($b -shl 16)
}
